# Applies the OktStayTypeCs CodeSystem update:
#  - bump Version to 0.2.0
#  - bump Date to 2025-11-12T16:10:30+00:00
#  - bump Count to 2
#  - add a new "korttidsopphold" concept row on the Concepts sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Cells.Item(3, 2).Value = "0.2.0"                           # Version
$meta.Cells.Item(8, 2).Value = "2025-11-12T16:10:30+00:00"       # Date

$meta.Cells.Item(22, 2).NumberFormat = "@"
$meta.Cells.Item(22, 2).Value = "2"                              # Count
$meta.Cells.Item(21, 2).Copy()
$meta.Cells.Item(22, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Concepts sheet ---------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# New row 3: keep text formatting consistent with the existing data row
$concepts.Cells.Item(3, 1).NumberFormat = "@"
$concepts.Cells.Item(3, 1).Value = "1"
$concepts.Cells.Item(3, 2).Value = "korttidsopphold"
$concepts.Cells.Item(3, 3).Value = "Korttidsopphold"
$concepts.Cells.Item(3, 4).Value = "Korttidsopphold i institusjon"

# Copy formatting (style/borders/wrap) from the existing data row onto the new one
$concepts.Range("A2:D2").Copy()
$concepts.Range("A3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
